$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "batsman" column (D), shifting
# batsman..sr from D:I to F:K, then fill the new ownTeam/oppTeam columns.
$ws.Range("D1:E1").EntireColumn.Insert()

$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("D2").Value = "Kings XI Punjab"
$ws.Range("E2").Value = "Chennai Super Kings"

# New match rows for KL Rahul (c) - numeric-looking columns must stay text
$ws.Range("G3:K4").NumberFormat = "@"

$ws.Range("A3").Value = " Abu Dhabi"
$ws.Range("B3").Value = " October 10 2020"
$ws.Range("C3").Value = "KKR won by 2 runs"
$ws.Range("D3").Value = "Kings XI Punjab"
$ws.Range("E3").Value = "Kolkata Knight Riders"
$ws.Range("F3").Value = "KL Rahul (c)"
$ws.Range("G3").Value = "74"
$ws.Range("H3").Value = "58"
$ws.Range("I3").Value = "6"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "127.58"

$ws.Range("A4").Value = " Dubai (DSC)"
$ws.Range("B4").Value = " October 08 2020"
$ws.Range("C4").Value = "Sunrisers won by 69 runs"
$ws.Range("D4").Value = "Kings XI Punjab"
$ws.Range("E4").Value = "Sunrisers Hyderabad"
$ws.Range("F4").Value = "KL Rahul (c)"
$ws.Range("G4").Value = "11"
$ws.Range("H4").Value = "16"
$ws.Range("I4").Value = "0"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "68.75"
